# Revert capacity chart to show kilowatts on the y-axis.
#
# 1. Worksheet data: the "Solar" column (E) values were stored in Watts;
#    convert the three non-zero entries to kilowatts (divide by 1000).
# 2. Number format used by the data grid (B:G) should show one decimal
#    place now that the values are fractional kilowatts.
# 3. The value (y) axis title changes from "Watts" to "Kilowatts (kW)".
# 4. The value (y) axis tick-label number format drops the ">=1000 ... K"
#    abbreviation and just shows "#,##0".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Convert the Solar (E) column capacity figures from W to kW -------
$ws.Range("E22").Value = 8.8
$ws.Range("E23").Value = 15.2
$ws.Range("E26").Value = 4.324

# --- 2. Show one decimal place across the numeric data grid --------------
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# --- 3 & 4. Update the chart's value-axis title and number format --------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$valAxis = $chart.Axes(2)
$valAxis.AxisTitle.Text = "Kilowatts (kW)"
$valAxis.TickLabels.NumberFormat = "#,##0"
